$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Session" to "Checklist"
$ws.Name = "Checklist"

# Insert a new row above the current row 2 (shifts old rows 2,3 down to 3,4)
$ws.Rows("2:2").Insert()

# Populate the newly inserted row 2 with the new scan/selection log entry.
# The leading apostrophe forces the Student ID to be stored as text (matching
# the text-typed "numeric looking" IDs already used in the other rows),
# exactly like typing '191153 into the cell in Excel.
$ws.Range("A2").Value = "'191153"
$ws.Range("B2").Value = "Pediatrics"
$ws.Range("C2").Value = "25/08/2025"
$ws.Range("D2").Value = "14:43:51"
$ws.Range("E2").Value = "Selection"
$ws.Range("F2").Value = "admin@admin.com"

# Update the "Type" column (E) for the remaining (now shifted down) rows,
# from "Scan" to "Selection"
$ws.Range("E3").Value = "Selection"
$ws.Range("E4").Value = "Selection"
